$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; this shifts the former rows 5 and 6
# (and their cell values/formatting) down to rows 6 and 7.
$ws.Rows.Item(5).Insert()

# Remove all existing hyperlinks on the sheet -- Insert() does not relocate
# the <hyperlinks> anchors, so we rebuild them from scratch below.
$ws.Hyperlinks.Delete()

# --- Write the refreshed scrape timestamp + row data ----------------------
# Row 2
$ws.Range('A2').Value = '2025-11-24 12:37:45'
$ws.Range('B2').Value = '【Python/AI/GAS 開発者・PM向け】「業務委託・再委託」の経験に関する30分インタビュー'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('G2').Value = 530
$ws.Range('H2').Value = '🔥AI,Python ◆開発'

# Row 3
$ws.Range('A3').Value = '2025-11-24 12:37:45'
$ws.Range('B3').Value = '【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('G3').Value = 158
$ws.Range('H3').Value = '◆自動化,スクレイピング ◇管理'

# Row 4
$ws.Range('A4').Value = '2025-11-24 12:37:45'
$ws.Range('B4').Value = 'マッチングサイト開発エンジニア募集'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('G4').Value = 100
$ws.Range('H4').Value = '◆開発 ◇サイト'

# Row 5
$ws.Range('A5').Value = '2025-11-24 12:37:45'
$ws.Range('B5').Value = '【急募】掲示板サイト(爆サイ)自動書き込みソフト開発者募集'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('G5').Value = 93
$ws.Range('H5').Value = '◆開発 ◇サイト'

# Row 6
$ws.Range('A6').Value = '2025-11-24 12:37:45'
$ws.Range('B6').Value = '急募 限定公開 PR 限定公開の仕事'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('G6').Value = 25

# Row 7
$ws.Range('A7').Value = '2025-11-24 12:37:45'
$ws.Range('B7').Value = '【急募】貸別荘収支表自動集計システム構築の依頼'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('G7').Value = 25

# Row 8
$ws.Range('A8').Value = '2025-11-24 12:37:45'
$ws.Range('B8').Value = '【急募】プログラム修正依頼!スキルを活かしてみませんか?'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('G8').Value = 13

# Row 9
$ws.Range('A9').Value = '2025-11-24 12:37:45'
$ws.Range('B9').Value = '【Amazon出品・Excel】ブラウズノード設定/フラットファイル検証に詳しい方を募集'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('G9').Value = 13

# --- Rebuild hyperlinks on column F in row order, matching the Hyperlink style ---
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5434693')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5440052')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5440077')
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5439484')
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5440230')
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5440042')
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5440002')
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5440204')
$ws.Range('F2:F9').Style = 'Hyperlink'
